$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.561.36'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.924.01'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4822'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08210'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.86'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.133'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').Value = '1.890.38'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.298'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06862'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001041'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.011'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = '29.553.08'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.669'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.52%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '2.126.96'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.411'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.093'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.013'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.617'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.566'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.379'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('E36').Value = '  +4.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02289'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.193'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5965'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.867'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1848'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.444'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.280'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07542'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5558'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.980'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '119.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('E51').Value = '  +0.80%  '
